$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(-0.37427986513949918, 0.3734497924577056),
    @(-0.22621075061616835, 0.22444117847417289),
    @(-0.12149629596769174, 0.12113742226844693),
    @(-0.1091374224604067, 0.10884019234868347),
    @(-0.10284019312983794, 0.10226753337479266),
    @(-0.031559878769990224, 0.031546856273856339),
    @(-0.011546857210262829, 0.011542204453757421),
    @(0.008457794609014968, -0.0084705181510855709),
    @(0.014470517354575385, -0.014496111830158753),
    @(0.020496111036521825, -0.020496073887755983),
    @(0.024996073109274874, -0.025015137187075709),
    @(0.031015136395133425, -0.031151671759050537),
    @(0.037151670975660522, -0.037221532646422695),
    @(-0.027078993459297429, 0.027049185276412757),
    @(-0.021049186058668568, 0.021025770668304133),
    @(-0.015025771452821246, 0.015003776771833177),
    @(-0.0090037775593501124, 0.008999999182050189),
    @(-0.13723904262895559, 0.13705100414413707),
    @(-0.12805100492507471, 0.12649792165418816),
    @(-0.053742647176562741, 0.053667717603115506),
    @(-0.044667718423431868, 0.044576515953222096),
    @(-0.09394488330373818, 0.093633257023135741),
    @(-0.084633257827499975, 0.084126382354111406),
    @(-0.042126383499205389, 0.041999998848523568),
    @(-0.064707655401431197, 0.064658410989352433),
    @(-0.058658411787003928, 0.058599057969445312),
    @(-0.052599058769597917, 0.052407208031094221),
    @(-0.046407208841536374, 0.046293224728741755),
    @(-0.034293225605345867, 0.034253108266934618),
    @(-0.014253109225712546, 0.013929974548507307),
    @(-0.052409719950448519, 0.052355339256028088),
    @(-0.031355340230041051, 0.031331791600192282)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
